$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row for 2022-Q4 at the
#    top of the data (row 2), pushing the existing quarters down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 35
$summary.Range("D2").Value = 9.17

# The newly inserted row inherited formatting from the insert operation;
# realign it with the look of the rest of the table: column A keeps the
# bold/centered/bordered style used by the other index cells, while B:D
# stay unstyled like the corresponding cells in the other rows.
$summary.Range("B2:D2").ClearFormats()
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Insert a brand-new "2022-Q4" sheet (fund holdings detail) right before
#    the existing "2022-Q3" sheet. We duplicate "2022-Q3" so the new sheet
#    starts out with identical layout/styling, then overwrite its contents.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.ActiveSheet
$q4.Name = "2022-Q4"

# The source sheet has 43 data rows (rows 2-44); the Q4 sheet only needs 35
# data rows (rows 2-36), so drop the extra trailing rows entirely.
$q4.Range("A37:H44").EntireRow.Delete()

$q4Data = @(
    @(0, "016464", "兴证全球合瑞混合A", "26.69", "85.54", "4.43", "1.1824", 4),
    @(1, "206002", "鹏华精选成长混合A", "31.14", "92.71", "3.59", "1.1179", 9),
    @(2, "010490", "鹏华高质量增长混合A", "13.18", "94.62", "8.38", "1.1045", 4),
    @(3, "398051", "中海环保新能源混合", "19.89", "77.11", "5.16", "1.0263", 3),
    @(4, "398021", "中海能源策略混合", "18.04", "90.92", "5.22", "0.9417", 2),
    @(5, "016465", "兴证全球合瑞混合C", "13.89", "85.54", "4.43", "0.6153", 4),
    @(6, "009984", "鹏华启航混合", "11.06", "90.05", "4.27", "0.4723", 5),
    @(7, "017732", "鹏华核心优势混合C", "5.55", "88.10", "7.37", "0.4090", 1),
    @(8, "011956", "鹏华新能源精选混合A", "6.74", "79.15", "5.82", "0.3923", 1),
    @(9, "011957", "鹏华新能源精选混合C", "5.72", "79.15", "5.82", "0.3329", 1),
    @(10, "009023", "鹏华稳健回报混合A", "4.33", "94.26", "7.12", "0.3083", 5),
    @(11, "011460", "鹏华创新成长混合A", "10.08", "78.91", "2.79", "0.2812", 6),
    @(12, "398061", "中海消费混合", "3.80", "90.62", "4.04", "0.1535", 8),
    @(13, "008811", "鹏华科技创新混合", "2.83", "89.03", "5.02", "0.1421", 3),
    @(14, "004986", "鹏华策略回报灵活配置混合", "3.00", "85.60", "3.78", "0.1134", 3),
    @(15, "000431", "鹏华品牌传承混合", "4.10", "86.27", "2.57", "0.1054", 10),
    @(16, "206012", "鹏华价值精选股票", "2.53", "87.28", "3.92", "0.0992", 4),
    @(17, "952035", "国泰君安君得诚混合", "2.21", "85.83", "4.03", "0.0891", 5),
    @(18, "016562", "鹏华精选成长混合C", "2.10", "92.71", "3.59", "0.0754", 9),
    @(19, "673110", "西部利得新润灵活配置混合A", "1.42", "80.39", "3.00", "0.0426", 9),
    @(20, "010491", "鹏华高质量增长混合C", "0.45", "94.62", "8.38", "0.0377", 4),
    @(21, "000166", "中海信息产业精选混合", "0.72", "83.53", "4.03", "0.0290", 2),
    @(22, "006526", "鹏华优选回报灵活配置混合A", "0.94", "89.78", "3.06", "0.0288", 10),
    @(23, "005175", "国寿安保消费新蓝海灵活配置混合", "0.70", "89.58", "2.91", "0.0204", 9),
    @(24, "970113", "兴证资管金麒麟兴睿优选一年持有期混合B", "0.67", "84.89", "2.93", "0.0196", 10),
    @(25, "004301", "国寿安保稳信混合A", "1.48", "22.04", "0.84", "0.0124", 10),
    @(26, "011461", "鹏华创新成长混合C", "0.35", "78.91", "2.79", "0.0098", 6),
    @(27, "017511", "鹏华稳健回报混合C", "0.05", "94.26", "7.12", "0.0036", 5),
    @(28, "012997", "鹏华优选回报灵活配置混合C", "0.04", "89.78", "3.06", "0.0012", 10),
    @(29, "004302", "国寿安保稳信混合C", "0.01", "22.04", "0.84", "0.0001", 10),
    @(30, "015356", "西部利得新润灵活配置混合C", "0.00", "80.39", "3.00", 0, 9),
    @(31, "970112", "兴证资管金麒麟兴睿优选一年持有期混合A", "0.00", "84.89", "2.93", 0, 10),
    @(32, "970114", "兴证资管金麒麟兴睿优选一年持有期混合C", "0.00", "84.89", "2.93", 0, 10),
    @(33, "015406", "国寿安保稳信混合E", "0.00", "22.04", "0.84", 0, 10),
    @(34, "006976", "鹏华核心优势混合A", "0.00", "88.10", "7.37", 0, 1)
)

$r = 2
foreach ($row in $q4Data) {
    # Columns that carry numeric-looking text (fund code, scale, position,
    # ratio, market value) must stay text so things like leading zeros in
    # fund codes are preserved instead of being coerced into numbers.
    $q4.Range("B" + $r + ":G" + $r).NumberFormat = "@"

    $q4.Range("A" + $r).Value = $row[0]
    $q4.Range("B" + $r).Value = $row[1]
    $q4.Range("C" + $r).Value = $row[2]
    $q4.Range("D" + $r).Value = $row[3]
    $q4.Range("E" + $r).Value = $row[4]
    $q4.Range("F" + $r).Value = $row[5]

    $gVal = $row[6]
    if ($gVal -eq 0) {
        $q4.Range("G" + $r).NumberFormat = "General"
        $q4.Range("G" + $r).Value = 0
    } else {
        $q4.Range("G" + $r).Value = $gVal
    }

    $q4.Range("H" + $r).Value = $row[7]

    $r = $r + 1
}
